$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'266.87"
$ws.Range("F2").Value = "'15-12-2022"
$ws.Range("G2").Value = "'0"

$ws.Range("D3").Value = "'22.71"
$ws.Range("F3").Value = "'15-12-2022"
$ws.Range("G3").Value = "'0"

$ws.Range("D4").Value = "'6.306"
$ws.Range("F4").Value = "'15-12-2022"
$ws.Range("G4").Value = "'0"

$ws.Range("D5").Value = "'0.06196"
$ws.Range("F5").Value = "'15-12-2022"
$ws.Range("G5").Value = "'0"

$ws.Range("D6").Value = "'3.587"
$ws.Range("F6").Value = "'15-12-2022"
$ws.Range("G6").Value = "'0"

$ws.Range("D7").Value = "'6.692"
$ws.Range("F7").Value = "'15-12-2022"
$ws.Range("G7").Value = "'0"

$ws.Range("D8").Value = "'1.366"
$ws.Range("F8").Value = "'15-12-2022"
$ws.Range("G8").Value = "'0"

$ws.Range("D9").Value = "'0.8419"
$ws.Range("F9").Value = "'15-12-2022"
$ws.Range("G9").Value = "'0"

$ws.Range("D10").Value = "'0.01361"
$ws.Range("F10").Value = "'15-12-2022"
$ws.Range("G10").Value = "'0"

$ws.Range("D11").Value = "'0.1604"
$ws.Range("F11").Value = "'15-12-2022"
$ws.Range("G11").Value = "'0"

$ws.Range("D12").Value = "'0.08259"
$ws.Range("F12").Value = "'15-12-2022"
$ws.Range("G12").Value = "'0"

$ws.Range("D13").Value = "'0.03422"
$ws.Range("F13").Value = "'15-12-2022"
$ws.Range("G13").Value = "'0"

$ws.Range("D14").Value = "'0.03214"
$ws.Range("F14").Value = "'15-12-2022"
$ws.Range("G14").Value = "'0"

$ws.Range("D15").Value = "'0.09263"
$ws.Range("F15").Value = "'15-12-2022"
$ws.Range("G15").Value = "'0"

$ws.Range("D16").Value = "'3.923"
$ws.Range("F16").Value = "'15-12-2022"
$ws.Range("G16").Value = "'0"

$ws.Range("D17").Value = "'0.001707"
$ws.Range("F17").Value = "'15-12-2022"
$ws.Range("G17").Value = "'0"

$ws.Range("D18").Value = "'0.04863"
$ws.Range("F18").Value = "'15-12-2022"
$ws.Range("G18").Value = "'0"

$ws.Range("D19").Value = "'0.006260"
$ws.Range("F19").Value = "'15-12-2022"
$ws.Range("G19").Value = "'0"

$ws.Range("D20").Value = "'0.005366"
$ws.Range("E20").Value = "'19HotbitTokenHTBWorstin24h"
$ws.Range("F20").Value = "'15-12-2022"
$ws.Range("G20").Value = "'0"

$ws.Range("F21").Value = "'15-12-2022"
$ws.Range("G21").Value = "'0"

$ws.Range("F22").Value = "'15-12-2022"
$ws.Range("G22").Value = "'0"

$ws.Range("D23").Value = "'3.768"
$ws.Range("F23").Value = "'15-12-2022"
$ws.Range("G23").Value = "'0"

$ws.Range("D24").Value = "'2.322"
$ws.Range("F24").Value = "'15-12-2022"
$ws.Range("G24").Value = "'0"

$ws.Range("F25").Value = "'15-12-2022"
$ws.Range("G25").Value = "'0"

$ws.Range("F26").Value = "'15-12-2022"
$ws.Range("G26").Value = "'0"

$ws.Range("F27").Value = "'15-12-2022"
$ws.Range("G27").Value = "'0"

$ws.Range("F28").Value = "'15-12-2022"
$ws.Range("G28").Value = "'0"

$ws.Range("F29").Value = "'15-12-2022"
$ws.Range("G29").Value = "'0"

$ws.Range("F30").Value = "'15-12-2022"
$ws.Range("G30").Value = "'0"

$ws.Range("F31").Value = "'15-12-2022"
$ws.Range("G31").Value = "'0"

$ws.Range("F32").Value = "'15-12-2022"
$ws.Range("G32").Value = "'0"

$ws.Range("F33").Value = "'15-12-2022"
$ws.Range("G33").Value = "'0"

$ws.Range("F34").Value = "'15-12-2022"
$ws.Range("G34").Value = "'0"

$ws.Range("F35").Value = "'15-12-2022"
$ws.Range("G35").Value = "'0"

$ws.Range("F36").Value = "'15-12-2022"
$ws.Range("G36").Value = "'0"

$ws.Range("F37").Value = "'15-12-2022"
$ws.Range("G37").Value = "'0"

$ws.Range("F38").Value = "'15-12-2022"
$ws.Range("G38").Value = "'0"

$ws.Range("F39").Value = "'15-12-2022"
$ws.Range("G39").Value = "'0"

$ws.Range("D40").Value = "'0.04655"
$ws.Range("F40").Value = "'15-12-2022"
$ws.Range("G40").Value = "'0"

$ws.Range("D41").Value = "'0.006969"
$ws.Range("F41").Value = "'15-12-2022"
$ws.Range("G41").Value = "'0"

$ws.Range("D42").Value = "'0.1153"
$ws.Range("F42").Value = "'15-12-2022"
$ws.Range("G42").Value = "'0"

$ws.Range("D43").Value = "'0.003203"
$ws.Range("F43").Value = "'15-12-2022"
$ws.Range("G43").Value = "'0"

$ws.Range("D44").Value = "'0.01111"
$ws.Range("F44").Value = "'15-12-2022"
$ws.Range("G44").Value = "'0"

$ws.Range("D45").Value = "'0.00006234"
$ws.Range("F45").Value = "'15-12-2022"
$ws.Range("G45").Value = "'0"

$ws.Range("F46").Value = "'15-12-2022"
$ws.Range("G46").Value = "'0"

$ws.Range("D47").Value = "'0.7890"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOIN"
$ws.Range("F47").Value = "'15-12-2022"
$ws.Range("G47").Value = "'0"

$ws.Range("D48").Value = "'0.1656"
$ws.Range("F48").Value = "'15-12-2022"
$ws.Range("G48").Value = "'0"

$ws.Range("D49").Value = "'0.00001400"
$ws.Range("F49").Value = "'15-12-2022"
$ws.Range("G49").Value = "'0"

$ws.Range("F50").Value = "'15-12-2022"
$ws.Range("G50").Value = "'0"

$ws.Range("F51").Value = "'15-12-2022"
$ws.Range("G51").Value = "'0"
